$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Shapes.Item("Connettore 1 129").Delete()
